$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (the blank rows), shifting rows 5/6 up to become 3/4
$ws.Rows("3:4").Delete()

# Delete columns A:D, shifting E:G left to become A:C
$ws.Columns("A:D").Delete()
